# 4th Test case (Create NSC)
# Adds a "Create New NSC" mini test-case table in columns K:N of rows 20-23
# on Sheet1, mirroring the existing "Create New INC" table (F8:I11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

# --- Values -----------------------------------------------------------
# Header row (row 20): Test Case Name | NSC Code | NSC Title |
$ws.Range("K20").Value = "Test Case Name"
$ws.Range("L20").Value = "NSC Code"
$ws.Range("M20").Value = "NSC Title"
$ws.Range("N20").Value = ""

# Data row (row 21): Create New NSC | 2229 | NSC 2229 Title Test |
$ws.Range("K21").Value = "Create New NSC"
$ws.Range("L21").Value = "2229"
$ws.Range("M21").Value = "NSC 2229 Title Test"
$ws.Range("N21").Value = ""

# Row 22 / 23 leftover helper cells
$ws.Range("K22").Value = ""
$ws.Range("L22").Value = "INC Code"
$ws.Range("M22").Value = ""
$ws.Range("N22").Value = ""

$ws.Range("K23").Value = ""
$ws.Range("L23").Value = "Test0001"
$ws.Range("M23").Value = ""
$ws.Range("N23").Value = ""

# --- Formatting (copy from the matching existing styled cells) --------
# NOTE: this host's PasteSpecial only honours the first area of a
# multi-area range, so every destination cell is pasted individually.
foreach ($dest in @("K20", "L20", "L22")) {
    $ws.Range("F10").Copy() | Out-Null
    $ws.Range($dest).PasteSpecial($xlPasteFormats) | Out-Null
}

foreach ($dest in @("M20", "N20")) {
    $ws.Range("H10").Copy() | Out-Null
    $ws.Range($dest).PasteSpecial($xlPasteFormats) | Out-Null
}

$ws.Range("F11").Copy() | Out-Null
$ws.Range("K21").PasteSpecial($xlPasteFormats) | Out-Null

foreach ($dest in @("L21", "L23")) {
    $ws.Range("G8").Copy() | Out-Null
    $ws.Range($dest).PasteSpecial($xlPasteFormats) | Out-Null
}

foreach ($dest in @("M21", "N21")) {
    $ws.Range("H11").Copy() | Out-Null
    $ws.Range($dest).PasteSpecial($xlPasteFormats) | Out-Null
}

foreach ($dest in @("K22", "K23")) {
    $ws.Range("F3").Copy() | Out-Null
    $ws.Range($dest).PasteSpecial($xlPasteFormats) | Out-Null
}

$ws.CutCopyMode = $false

# --- Row heights --------------------------------------------------------
$ws.Rows.Item(21).RowHeight = 12.75
$ws.Rows.Item(23).RowHeight = 12.75

# --- View / selection ----------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 7
$ws.Range("L21").Select() | Out-Null
